# Applies the "Automatic update of files" change to the Avverkningsanmälningar sheet:
#  - Column C ("Förändrad") bumps from 46065 to 46066 for every data row (2-9).
#  - Rows 4 and 9 swap their Beteckning (A), Datum (B) and Area (G) values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

foreach ($r in 2..9) {
    $ws.Cells.Item($r, 3).Value = 46066
}

$ws.Cells.Item(4, 1).Value = "A 25617-2024"
$ws.Cells.Item(4, 2).Value = 45463
$ws.Cells.Item(4, 7).Value = 2.3

$ws.Cells.Item(9, 1).Value = "A 17908-2021"
$ws.Cells.Item(9, 2).Value = 44301
$ws.Cells.Item(9, 7).Value = 0.9
